$wb = $excel.ActiveWorkbook

# --- watchlist sheet: re-sorted by RSI ascending with refreshed quotes ---
$ws1 = $wb.Worksheets.Item("watchlist")
$ws1.Range("A2").Value = "CSCO"
$ws1.Range("B2").Value = 41.87
$ws1.Range("C2").Value = -0.262
$ws1.Range("D2").Value = 18.8356164383562
$ws1.Range("A3").Value = "XOM"
$ws1.Range("B3").Value = 41.96
$ws1.Range("C3").Value = -1.1077
$ws1.Range("D3").Value = 38.23529411764709
$ws1.Range("A4").Value = "BA"
$ws1.Range("B4").Value = 169.27
$ws1.Range("C4").Value = -0.5639
$ws1.Range("D4").Value = 42.15471343564047
$ws1.Range("A5").Value = "TRV"
$ws1.Range("B5").Value = 113.2
$ws1.Range("C5").Value = -1.5053
$ws1.Range("D5").Value = 43.75830013280214
$ws1.Range("A6").Value = "MSFT"
$ws1.Range("B6").Value = 209.7
$ws1.Range("C6").Value = -0.8464
$ws1.Range("D6").Value = 44.05722670579603
$ws1.Range("A7").Value = "PFE"
$ws1.Range("B7").Value = 38.26
$ws1.Range("C7").Value = -0.2607
$ws1.Range("D7").Value = 45.73991031390123
$ws1.Range("A8").Value = "GS"
$ws1.Range("B8").Value = 203.02
$ws1.Range("C8").Value = 0.8494
$ws1.Range("D8").Value = 46.37330754352031
$ws1.Range("A9").Value = "INTC"
$ws1.Range("B9").Value = 48.33
$ws1.Range("C9").Value = -0.6578
$ws1.Range("D9").Value = 47.58364312267653
$ws1.Range("A10").Value = "WBA"
$ws1.Range("B10").Value = 40.25
$ws1.Range("C10").Value = -1.2512
$ws1.Range("D10").Value = 49.17218543046356
$ws1.Range("A11").Value = "IBM"
$ws1.Range("B11").Value = 123.84
$ws1.Range("C11").Value = -0.8646
$ws1.Range("D11").Value = 50.06257822277852
$ws1.Range("A12").Value = "DIS"
$ws1.Range("B12").Value = 127.77
$ws1.Range("C12").Value = -0.892
$ws1.Range("D12").Value = 50.65252854812396
$ws1.Range("A13").Value = "CVX"
$ws1.Range("B13").Value = 86.39
$ws1.Range("C13").Value = -1.415
$ws1.Range("D13").Value = 52.16952573158428
$ws1.Range("A14").Value = "KO"
$ws1.Range("B14").Value = 47.37
$ws1.Range("C14").Value = -2.1685
$ws1.Range("D14").Value = 52.50836120401335
$ws1.Range("A15").Value = "AXP"
$ws1.Range("B15").Value = 96.89
$ws1.Range("C15").Value = 0.2276
$ws1.Range("D15").Value = 55.10899182561308
$ws1.Range("A16").Value = "UNH"
$ws1.Range("B16").Value = 315.4
$ws1.Range("C16").Value = -0.4451
$ws1.Range("D16").Value = 55.20796302879479
$ws1.Range("A17").Value = "JPM"
$ws1.Range("B17").Value = 98.55
$ws1.Range("C17").Value = 0.2339
$ws1.Range("D17").Value = 55.39452495974237
$ws1.Range("A18").Value = "JNJ"
$ws1.Range("B18").Value = 150.39
$ws1.Range("C18").Value = 0.1999
$ws1.Range("D18").Value = 59.68841285296976
$ws1.Range("A19").Value = "WMT"
$ws1.Range("B19").Value = 132.41
$ws1.Range("C19").Value = -1.7074
$ws1.Range("D19").Value = 61.36528685548289
$ws1.Range("A20").Value = "CAT"
$ws1.Range("B20").Value = 138.02
$ws1.Range("C20").Value = -0.2529
$ws1.Range("D20").Value = 61.4575507137491
$ws1.Range("A21").Value = "RTX"
$ws1.Range("B21").Value = 61.64
$ws1.Range("C21").Value = -0.7407
$ws1.Range("D21").Value = 63.87009472259805
$ws1.Range("A22").Value = "PG"
$ws1.Range("B22").Value = 135.77
$ws1.Range("C22").Value = -0.5421
$ws1.Range("D22").Value = 63.91875746714468
$ws1.Range("A23").Value = "AAPL"
$ws1.Range("B23").Value = 462.83
$ws1.Range("C23").Value = 0.1255
$ws1.Range("D23").Value = 65.77135915778103
$ws1.Range("A24").Value = "V"
$ws1.Range("B24").Value = 200.99
$ws1.Range("C24").Value = 0.9949
$ws1.Range("D24").Value = 67.95847750865045
$ws1.Range("A25").Value = "VZ"
$ws1.Range("B25").Value = 59.05
$ws1.Range("C25").Value = -0.2197
$ws1.Range("D25").Value = 75.94501718213053
$ws1.Range("A26").Value = "HD"
$ws1.Range("B26").Value = 282.86
$ws1.Range("C26").Value = -0.7509
$ws1.Range("D26").Value = 77.11565585331451
$ws1.Range("A27").Value = "DOW"
$ws1.Range("B27").Value = 44.64
$ws1.Range("C27").Value = -0.1566
$ws1.Range("D27").Value = 78.49462365591407
$ws1.Range("A28").Value = "MRK"
$ws1.Range("B28").Value = 85.03
$ws1.Range("C28").Value = 0.5558
$ws1.Range("D28").Value = 81.56424581005577
$ws1.Range("A29").Value = "MMM"
$ws1.Range("B29").Value = 163.97
$ws1.Range("C29").Value = -0.2494
$ws1.Range("D29").Value = 83.46273291925469
$ws1.Range("A30").Value = "NKE"
$ws1.Range("B30").Value = 108.39
$ws1.Range("C30").Value = 1.3275
$ws1.Range("D30").Value = 84.71575023299171
$ws1.Range("A31").Value = "MCD"
$ws1.Range("B31").Value = 209.51
$ws1.Range("C31").Value = -0.3851
$ws1.Range("D31").Value = 89.9454403741232

# --- stocks sheet: refreshed current price / value / performance / rsi ---
$ws2 = $wb.Worksheets.Item("stocks")
$ws2.Range("C2").Value = 41.96
$ws2.Range("E2").Value = 503.52
$ws2.Range("F2").Value = -8.7429
$ws2.Range("G2").Value = 38.23529411764709
$ws2.Range("C3").Value = 48.33
$ws2.Range("E3").Value = 869.9399999999999
$ws2.Range("F3").Value = -4.4673
$ws2.Range("G3").Value = 47.58364312267653
$ws2.Range("C4").Value = 150.39
$ws2.Range("E4").Value = 601.56
$ws2.Range("F4").Value = 2.4176
$ws2.Range("G4").Value = 59.68841285296976
$ws2.Range("C5").Value = 113.2
$ws2.Range("E5").Value = 452.8
$ws2.Range("F5").Value = 0.0884
$ws2.Range("G5").Value = 43.75830013280214
$ws2.Range("C6").Value = 38.26
$ws2.Range("E6").Value = 535.64
$ws2.Range("F6").Value = 1.2437
$ws2.Range("G6").Value = 45.73991031390123
$ws2.Range("C7").Value = 41.87
$ws2.Range("E7").Value = 586.18
$ws2.Range("F7").Value = -1.4824
$ws2.Range("G7").Value = 18.8356164383562

# --- portfolio sheet: stock value + total refreshed ---
$ws3 = $wb.Worksheets.Item("portfolio")
$ws3.Range("B3").Value = 3549.64
$ws3.Range("B4").Value = 10981.27

# --- summary sheet: append 08-18 and 08-19 snapshots ---
$ws5 = $wb.Worksheets.Item("summary")
$ws5.Range("A65:D65").Copy($ws5.Range("A66:D66"))
$ws5.Range("A65:D65").Copy($ws5.Range("A67:D67"))
$ws5.Range("A66").Value = "18/08/2020 17:22:15"
$ws5.Range("B66").Value = 7431.629999999999
$ws5.Range("C66").Value = 3569.699999999999
$ws5.Range("D66").Value = 11001.33
$ws5.Range("A67").Value = "19/08/2020 17:04:31"
$ws5.Range("B67").Value = 7431.629999999999
$ws5.Range("C67").Value = 3549.64
$ws5.Range("D67").Value = 10981.27

Write-Output "edit complete"
